$wb = $excel.ActiveWorkbook

# --- Step 1: rename "cumcontrol" to "cumcontrol1" ---
$ws32 = $wb.Worksheets.Item("cumcontrol")
$ws32.Name = "cumcontrol1"

# --- Step 2: duplicate cumcontrol1 to create cumcontrol2, placed right after it ---
# (this preserves the original cumcontrol text/formatting as the starting point,
#  which is then edited below into the "dickpic"-style DELAY/SYNC/EDGE copy)
$ws32.Copy($null, $ws32)
$ws33 = $wb.Worksheets.Item($ws32.Index + 1)
$ws33.Name = "cumcontrol2"

# --- References to the other sheets touched by this edit ---
$ws1 = $wb.Worksheets.Item("JessicaFPJourney")
$ws5 = $wb.Worksheets.Item("ReEngagement")

# ==== sheet1 (JessicaFPJourney) edits ====
$ws1.Range("B2").Value = 'I seriously needed that mi amor. you''re different from everyone else on here, that was real'
$ws1.Range("B4").Value = 'I''m about to finish... stay with me'
$ws1.Range("B5").Value = 'hold on hold on'
$ws1.Range("B6").Value = 'dios mio I wanna cum with you... I''m so close'
$ws1.Range("B10").Value = 'hold on a sec'
$ws1.Range("B11").Value = 'I can''t control myself anymore'
$ws1.Range("B12").Value = 'imagine me pinning you down and not letting you move while I ride you amor'
$ws1.Range("B19").Value = 'ay dios... talking to you right now is dangerous'
$ws1.Range("B22").Value = 'hold on a sec'
$ws1.Range("B30").Value = 'I just came back from tango class and I''m still buzzing... everything feels so sensual right now and I keep thinking about you'
$ws1.Range("C30").Value = 'THE PIVOT. Tango = physical/sensual bridge.'
$ws1.Range("B31").Value = 'I barely get free time with everything I have going on but honestly? I''d rather spend it talking to you than anyone else right now'
$ws1.Range("C31").Value = 'Ego boost. Hints at responsibilities (single mom) without details. Transition.'
$ws1.Range("B33").Value = 'nice! I''m from Tucuman, Argentina. gym instructor by day but sailing is really my passion... there''s nothing like being out on the water'

# ==== cumcontrol1 edits ====
$ws32.Range("B2").Value = 'patience... what''s coming is worth every second of waiting'
$ws32.Range("B3").Value = 'hold it for me... I have years of experience and this next one is my best work'
$ws32.Range("C3").Value = 'DELAY. Send PPV.'
$ws32.Range("B4").Value = 'I want to feel you let go while I do the same amor... watch this first'
$ws32.Range("C4").Value = 'SYNC variant. Send PPV.'
$ws32.Range("B5").Value = 'now we go together... I''ve been holding back too. open this'
$ws32.Range("C5").Value = 'SYNC. Send PPV.'
$ws32.Range("B6").Value = 'a man who can wait gets rewarded... trust me on that'
$ws32.Range("B7").Value = 'I can tell you''re close... not yet amor, I know what I''m doing'
$ws32.Range("C7").Value = 'CONTROL.'

# ==== cumcontrol2 edits ====
$ws33.Range("A2").Value = 'delay2'
$ws33.Range("B2").Value = 'save it for this last one amor, I promise you it''s going to be worth it'
$ws33.Range("C2").Value = 'DELAY variant.'
$ws33.Range("A3").Value = 'delay1'
$ws33.Range("B3").Value = 'one more for you before we''re done... this is the one I''m most proud of'
$ws33.Range("C3").Value = 'DELAY. Send PPV.'
$ws33.Range("A4").Value = 'sync2'
$ws33.Range("B4").Value = 'I''m ready when you are... but see this first'
$ws33.Range("C4").Value = 'SYNC variant.'
$ws33.Range("A5").Value = 'sync1'
$ws33.Range("B5").Value = 'okay amor... let''s both let go right now. open this'
$ws33.Range("C5").Value = 'SYNC. Send PPV.'
$ws33.Range("A6").Value = 'edge2'
$ws33.Range("B6").Value = 'not yet... a little more anticipation makes it so much better, trust me'
$ws33.Range("C6").Value = 'EDGE variant.'
$ws33.Range("A7").Value = 'edge1'
$ws33.Range("B7").Value = 'slow down for me... I know exactly when to let you go'
$ws33.Range("C7").Value = 'CONTROL.'

# ==== ReEngagement edit ====
$ws5.Range("B3").Value = 'just got back from the marina and can''t stop thinking about you... you free?'

Write-Host "Done. Sheet order:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Index ":" $s.Name
}
